$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("unit_file")
$ws.Rows("10:18").Delete()
$ws.Range("E13").Select()
